$wb = $excel.ActiveWorkbook

# This script applies refreshed market-price/profit figures (columns H-N)
# to specific rows across all 8 sheets, matching a scheduled data-refresh commit.

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 822.9
$ws.Range("I2").Value = 1313
$ws.Range("J2").Value = 87.75
$ws.Range("K2").Value = 1313
$ws.Range("L2").Value = 87.75
$ws.Range("M2").Value = -1200
$ws.Range("N2").Value = -313.75
$ws.Range("H9").Value = 243.14285
$ws.Range("I9").Value = 384
$ws.Range("J9").Value = 186.8
$ws.Range("K9").Value = 384
$ws.Range("L9").Value = 186.8
$ws.Range("M9").Value = -215
$ws.Range("N9").Value = -524.8
$ws.Range("H43").Value = 9500
$ws.Range("I43").Value = 13500
$ws.Range("J43").Value = 8928.571
$ws.Range("K43").Value = 13500
$ws.Range("L43").Value = 8928.571
$ws.Range("M43").Value = -13431
$ws.Range("N43").Value = -9066.571
$ws.Range("H46").Value = 0
$ws.Range("I46").Value = 0
$ws.Range("K46").Value = 0
$ws.Range("M46").ClearContents()
$ws.Range("H53").Value = 294.35
$ws.Range("I53").Value = 228.27272
$ws.Range("J53").Value = 375.1111
$ws.Range("K53").Value = 228.27272
$ws.Range("L53").Value = 375.1111
$ws.Range("M53").Value = 408.72728
$ws.Range("N53").Value = -1649.1111
$ws.Range("H60").Value = 0
$ws.Range("I60").Value = 0
$ws.Range("K60").Value = 0
$ws.Range("M60").ClearContents()
$ws.Range("H132").Value = 3915.7646
$ws.Range("I132").Value = 3183.5715
$ws.Range("J132").Value = 7332.6665
$ws.Range("K132").Value = 9550.7145
$ws.Range("L132").Value = 21997.9995
$ws.Range("M132").Value = -7020.7145
$ws.Range("N132").Value = -27057.9995
$ws.Range("H133").Value = 90397.60000000001
$ws.Range("J133").Value = 90397.60000000001
$ws.Range("L133").Value = 90397.60000000001
$ws.Range("N133").Value = -100517.6
$ws.Range("H134").Value = 97373.5
$ws.Range("J134").Value = 69999
$ws.Range("L134").Value = 69999
$ws.Range("N134").Value = -80139
$ws.Range("H135").Value = 1244.0741
$ws.Range("I135").Value = 1166.72
$ws.Range("K135").Value = 10500.48
$ws.Range("M135").Value = -7965.48
$ws.Range("H137").Value = 2435.077
$ws.Range("I137").Value = 2013.9474
$ws.Range("J137").Value = 3578.1428
$ws.Range("K137").Value = 6041.8422
$ws.Range("L137").Value = 10734.4284
$ws.Range("M137").Value = -3491.8422
$ws.Range("N137").Value = -15834.4284
$ws.Range("H138").Value = 2706.4482
$ws.Range("I138").Value = 1728.2
$ws.Range("J138").Value = 3754.5715
$ws.Range("K138").Value = 5184.6
$ws.Range("L138").Value = 11263.7145
$ws.Range("M138").Value = -44.60000000000036
$ws.Range("N138").Value = -21543.7145

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2700
$ws.Range("I2").Value = 2509.5454
$ws.Range("J2").Value = 3747.5
$ws.Range("K2").Value = 2509.5454
$ws.Range("L2").Value = 3747.5
$ws.Range("M2").Value = -2396.5454
$ws.Range("N2").Value = -3973.5
$ws.Range("H31").Value = 3338.1667
$ws.Range("I31").Value = 3338.1667
$ws.Range("J31").Value = 0
$ws.Range("K31").Value = 3338.1667
$ws.Range("L31").Value = 0
$ws.Range("M31").ClearContents()
$ws.Range("N31").Value = -3044.1667
$ws.Range("H32").Value = 19409.887
$ws.Range("I32").Value = 3542.9333
$ws.Range("K32").Value = 3542.9333
$ws.Range("M32").Value = -3255.9333
$ws.Range("H45").Value = 528401.9399999999
$ws.Range("I45").Value = 770664.3
$ws.Range("K45").Value = 770664.3
$ws.Range("M45").Value = -770287.3
$ws.Range("H61").Value = 1981.8298
$ws.Range("I61").Value = 1732.3572
$ws.Range("K61").Value = 1732.3572
$ws.Range("M61").Value = -1520.3572
$ws.Range("H80").Value = 19997.143
$ws.Range("J80").Value = 19997.143
$ws.Range("L80").Value = 19997.143
$ws.Range("N80").Value = -21993.143
$ws.Range("H83").Value = 19997.143
$ws.Range("J83").Value = 19997.143
$ws.Range("L83").Value = 59991.429
$ws.Range("N83").Value = -69975.429
$ws.Range("H110").Value = 2087.9
$ws.Range("I110").Value = 2087.9
$ws.Range("K110").Value = 2087.9
$ws.Range("M110").Value = -42.90000000000009
$ws.Range("H116").Value = 2700
$ws.Range("I116").Value = 2509.5454
$ws.Range("J116").Value = 3747.5
$ws.Range("K116").Value = 2509.5454
$ws.Range("L116").Value = 3747.5
$ws.Range("M116").Value = -215.5454
$ws.Range("N116").Value = -8335.5
$ws.Range("H122").Value = 3333.6667
$ws.Range("I122").Value = 3333.6667
$ws.Range("K122").Value = 10001.0001
$ws.Range("M122").Value = -7551.000100000001
$ws.Range("H132").Value = 8276.134
$ws.Range("I132").Value = 8733.393
$ws.Range("K132").Value = 26200.179
$ws.Range("M132").Value = -23670.179
$ws.Range("H136").Value = 1981.8298
$ws.Range("I136").Value = 1732.3572
$ws.Range("K136").Value = 5197.071599999999
$ws.Range("M136").Value = -2647.071599999999

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2700
$ws.Range("I3").Value = 2509.5454
$ws.Range("J3").Value = 3747.5
$ws.Range("K3").Value = 2509.5454
$ws.Range("L3").Value = 3747.5
$ws.Range("M3").Value = -2395.5454
$ws.Range("N3").Value = -3975.5
$ws.Range("H95").Value = 45311.75
$ws.Range("J95").Value = 45311.75
$ws.Range("L95").Value = 45311.75
$ws.Range("N95").Value = -50803.75
$ws.Range("H105").Value = 5563.1055
$ws.Range("I105").Value = 5649.3335
$ws.Range("K105").Value = 5649.3335
$ws.Range("M105").Value = -3902.3335
$ws.Range("H141").Value = 66948
$ws.Range("J141").Value = 69027.664
$ws.Range("L141").Value = 69027.664
$ws.Range("N141").Value = -79387.664

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1582.5652
$ws.Range("I58").Value = 1485.7693
$ws.Range("J58").Value = 1708.4
$ws.Range("K58").Value = 1485.7693
$ws.Range("L58").Value = 1708.4
$ws.Range("M58").Value = -1282.7693
$ws.Range("N58").Value = -2114.4
$ws.Range("H59").Value = 18074.309
$ws.Range("J59").Value = 17913.916
$ws.Range("L59").Value = 17913.916
$ws.Range("N59").Value = -20203.916
$ws.Range("H60").Value = 8383.846
$ws.Range("I60").Value = 2999.3333
$ws.Range("J60").Value = 9999.200000000001
$ws.Range("K60").Value = 2999.3333
$ws.Range("L60").Value = 9999.200000000001
$ws.Range("M60").Value = -2488.3333
$ws.Range("N60").Value = -11021.2
$ws.Range("H134").Value = 2384
$ws.Range("I134").Value = 2459.88
$ws.Range("K134").Value = 7379.64
$ws.Range("M134").Value = -4844.64
$ws.Range("H136").Value = 1582.5652
$ws.Range("I136").Value = 1485.7693
$ws.Range("J136").Value = 1708.4
$ws.Range("K136").Value = 4457.3079
$ws.Range("L136").Value = 5125.200000000001
$ws.Range("M136").Value = -1907.3079
$ws.Range("N136").Value = -10225.2
$ws.Range("H141").Value = 167648.8
$ws.Range("J141").Value = 167648.8
$ws.Range("L141").Value = 167648.8
$ws.Range("N141").Value = -178008.8

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 1616.6875
$ws.Range("I34").Value = 110.14286
$ws.Range("J34").Value = 2788.4443
$ws.Range("K34").Value = 330.42858
$ws.Range("L34").Value = 8365.332900000001
$ws.Range("M34").Value = -246.42858
$ws.Range("N34").Value = -8533.332900000001
$ws.Range("H68").Value = 1749.5
$ws.Range("J68").Value = 1999
$ws.Range("L68").Value = 5997
$ws.Range("N68").Value = -7619
$ws.Range("H70").Value = 4760
$ws.Range("I70").Value = 3924.6667
$ws.Range("J70").Value = 6013
$ws.Range("K70").Value = 11774.0001
$ws.Range("L70").Value = 18039
$ws.Range("M70").Value = -11459.0001
$ws.Range("N70").Value = -18669
$ws.Range("H71").Value = 1749.5
$ws.Range("J71").Value = 1999
$ws.Range("L71").Value = 17991
$ws.Range("N71").Value = -26103
$ws.Range("H73").Value = 4760
$ws.Range("I73").Value = 3924.6667
$ws.Range("J73").Value = 6013
$ws.Range("K73").Value = 11774.0001
$ws.Range("L73").Value = 18039
$ws.Range("M73").Value = -10682.0001
$ws.Range("N73").Value = -20223
$ws.Range("H80").Value = 2000
$ws.Range("I80").Value = 2000
$ws.Range("K80").Value = 6000
$ws.Range("M80").Value = -5064
$ws.Range("H83").Value = 2000
$ws.Range("I83").Value = 2000
$ws.Range("K83").Value = 18000
$ws.Range("M83").Value = -13320
$ws.Range("H113").Value = 1059.8
$ws.Range("I113").Value = 1466.3334
$ws.Range("J113").Value = 958.1667
$ws.Range("K113").Value = 4399.0002
$ws.Range("L113").Value = 2874.5001
$ws.Range("M113").Value = -2229.0002
$ws.Range("N113").Value = -7214.5001
$ws.Range("H131").Value = 12346.296
$ws.Range("J131").Value = 35698
$ws.Range("L131").Value = 107094
$ws.Range("N131").Value = -117174

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2465.9473
$ws.Range("I132").Value = 2248.3713
$ws.Range("J132").Value = 5004.3335
$ws.Range("K132").Value = 6745.113899999999
$ws.Range("L132").Value = 15013.0005
$ws.Range("M132").Value = -4215.113899999999
$ws.Range("N132").Value = -20073.0005

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 2100.25
$ws.Range("J46").Value = 901
$ws.Range("L46").Value = 901
$ws.Range("N46").Value = -1277
$ws.Range("H132").Value = 2777.585
$ws.Range("I132").Value = 2294.5122
$ws.Range("K132").Value = 6883.5366
$ws.Range("M132").Value = -4353.5366
$ws.Range("H136").Value = 4102.478
$ws.Range("I136").Value = 3072.0833
$ws.Range("K136").Value = 9216.249899999999
$ws.Range("M136").Value = -6666.249899999999

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 11701.412
$ws.Range("J126").Value = 24991
$ws.Range("L126").Value = 74973
$ws.Range("N126").Value = -79913
$ws.Range("H132").Value = 11025.6875
$ws.Range("I132").Value = 16501.4
$ws.Range("J132").Value = 1899.5
$ws.Range("K132").Value = 49504.2
$ws.Range("L132").Value = 5698.5
$ws.Range("M132").Value = -46974.2
$ws.Range("N132").Value = -10758.5
